# "clients" sheet (3rd sheet / index 3 in 1-based Worksheets collection).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Trim trailing spaces off a handful of name entries ---------------
$ws.Range("B2").Value = "Kataryna"
$ws.Range("C2").Value = "Anna"
$ws.Range("B3").Value = "Maria"
$ws.Range("B5").Value = "Aganemnon"

# --- Phone numbers (column D): were plain numbers, now grouped text ---
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("D2").Value = "123 123 123"
$ws.Range("D3").Value = "799 799 799"
$ws.Range("D4").Value = "888 888 999"
$ws.Range("D5").Value = "320 129 399"
$ws.Range("D6").Value = "444 232 421"

# --- Postal code that used to be stored as a bare number ---------------
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "12-222"

# --- "is company" column (column E): 0/number -> literal text FALSE ---
$ws.Range("E2:E6").NumberFormat = "@"
$ws.Range("E2").Value = "'FALSE"
$ws.Range("E3").Value = "'FALSE"
$ws.Range("E4").Value = "'FALSE"
$ws.Range("E5").Value = "'FALSE"
$ws.Range("E6").Value = "'FALSE"

# --- Selection moved from F15 to D3 ------------------------------------
$ws.Range("D3").Select()

# --- Page setup: paper size + portrait orientation ---------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
